$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row
$ws.Range("A1").Value = "Venture Order List"

# Header row
$ws.Range("A2").Value = "orderId"
$ws.Range("B2").Value = "partId"
$ws.Range("C2").Value = "projectName"
$ws.Range("D2").Value = "lastMaterialDate"
$ws.Range("E2").Value = "shipDate"
$ws.Range("F2").Value = "quantity"

# Data rows
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 3003
$ws.Range("C3").Value = "Test order 1"
$ws.Range("D3").Value = "15/4/2020  12:00PM"
$ws.Range("E3").Value = "30/4/2020  12:12PM"
$ws.Range("F3").Value = 999

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 3007
$ws.Range("C4").Value = "Test order 2"
$ws.Range("D4").Value = "12/3/2020  12:00PM"
$ws.Range("E4").Value = "30/3/2020  12:00PM"
$ws.Range("F4").Value = 999

$ws.Range("A5").Value = 13
$ws.Range("B5").Value = 3008
$ws.Range("C5").Value = "Test order 3"
$ws.Range("D5").Value = "12/3/2020  12:00PM"
$ws.Range("E5").Value = "16/12/2020  12:12PM"
$ws.Range("F5").Value = 999

$ws.Range("A6").Value = 14
$ws.Range("B6").Value = 3008
$ws.Range("C6").Value = "Test order 4"
$ws.Range("D6").Value = "1/1/2020  12:00PM"
$ws.Range("E6").Value = "30/03/2020  12:00PM"
$ws.Range("F6").Value = 999

# Apply the same bold formatting that B3:B6 needs (matches style index 2, same as header formatting)
$ws.Range("B3:B6").Font.Bold = $true

# Update selection to match the post-edit state
$ws.Activate()
$ws.Range("A7:XFD54").Select()
